# implement dsa / tornadoplot
#
# Updates the "price" (F) column of the cost/utilisation diagnostics sheet
# with the recomputed DSA / tornado-plot figures, drops the inherited
# "-webkit-standard" 14pt style from most of those cells (two cells per
# 8-row diagnosis block keep it), tightens the row height of the rows
# whose style changed from 19 to 17 points, and moves the sheet's
# viewport/selection to match where the author ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F-column (price) value + style updates -------------------------------
# Every 8-row diagnosis block (2-9, 10-17, 18-25, 26-33, 34-41) gets the same
# sequence of new values. The 1st and 7th row of each block keep the
# original "s=4" cell style; the rest fall back to the workbook's default
# (Normal) style.
$fColumnUpdates = @(
    @{ Row = 2; Value = 0; KeepStyle = $true },
    @{ Row = 3; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 4; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 5; Value = 39.020000000000003; KeepStyle = $false },
    @{ Row = 6; Value = 79.510000000000005; KeepStyle = $false },
    @{ Row = 7; Value = 40.47; KeepStyle = $false },
    @{ Row = 8; Value = 171.71; KeepStyle = $true },
    @{ Row = 9; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 10; Value = 0; KeepStyle = $true },
    @{ Row = 11; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 12; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 13; Value = 39.020000000000003; KeepStyle = $false },
    @{ Row = 14; Value = 79.510000000000005; KeepStyle = $false },
    @{ Row = 15; Value = 40.47; KeepStyle = $false },
    @{ Row = 16; Value = 171.71; KeepStyle = $true },
    @{ Row = 17; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 18; Value = 0; KeepStyle = $true },
    @{ Row = 19; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 20; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 21; Value = 39.020000000000003; KeepStyle = $false },
    @{ Row = 22; Value = 79.510000000000005; KeepStyle = $false },
    @{ Row = 23; Value = 40.47; KeepStyle = $false },
    @{ Row = 24; Value = 171.71; KeepStyle = $true },
    @{ Row = 25; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 26; Value = 0; KeepStyle = $true },
    @{ Row = 27; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 28; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 29; Value = 39.020000000000003; KeepStyle = $false },
    @{ Row = 30; Value = 79.510000000000005; KeepStyle = $false },
    @{ Row = 31; Value = 40.47; KeepStyle = $false },
    @{ Row = 32; Value = 171.71; KeepStyle = $true },
    @{ Row = 33; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 34; Value = 0; KeepStyle = $true },
    @{ Row = 35; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 36; Value = 40.016666666666602; KeepStyle = $false },
    @{ Row = 37; Value = 39.020000000000003; KeepStyle = $false },
    @{ Row = 38; Value = 79.510000000000005; KeepStyle = $false },
    @{ Row = 39; Value = 40.47; KeepStyle = $false },
    @{ Row = 40; Value = 171.71; KeepStyle = $true },
    @{ Row = 41; Value = 40.016666666666602; KeepStyle = $false }
)

foreach ($item in $fColumnUpdates) {
    $cell = $ws.Cells.Item($item.Row, 6)
    $cell.Value = $item.Value
    if (-not $item.KeepStyle) {
        $cell.Style = "Normal"
    }
}

# --- Row height tightening -------------------------------------------------
# The rows whose "price" cell lost its s=4 style and used to be 19pt tall
# (the short single-line header rows of each diagnosis block) shrink to 17pt.
$heightRows = @(3, 7, 11, 15, 19, 23, 27, 31, 35, 39)
foreach ($r in $heightRows) {
    $ws.Rows($r).RowHeight = 17
}

# --- Viewport / selection ---------------------------------------------------
# Move the working selection to where the author left off.
$ws.Range("A11").Select()
$ws.Range("J33").Select()
